$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Rarity column values to lowercase ("Epic"/"Rare"/"Uncommon" -> "epic"/"rare"/"uncommon")
$ws.Range("C2:C6").Value = "epic"
$ws.Range("C7:C11").Value = "rare"
$ws.Range("C12:C16").Value = "uncommon"

# Remove the "READ FIRST" instructions text box from the sheet
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# Move the active selection from B18 to E14
$ws.Range("E14").Select()
